# Update gh-pages to output generated at 456a3b4
# Refresh the scraped 哔哩哔哩 (bilibili) "想去人数" (interest counters) and
# one event cover-image URL across the "展览" (sheet 1) and "全部类型"
# (sheet 4) worksheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# ---- 展览 (sheet 1) ----
$ws1.Range("F5").Value = 15345
$ws1.Range("F8").Value = 689
$ws1.Range("F9").Value = 15332
$ws1.Range("F11").Value = 8911
$ws1.Range("F17").Value = 192
$ws1.Range("F31").Value = 45
$ws1.Range("I33").Value = "//i0.hdslb.com/bfs/openplatform/202409/52AMZyUi1727059410434.jpeg"
$ws1.Range("F34").Value = 295
$ws1.Range("F37").Value = 5452

# ---- 全部类型 (sheet 4) ----
$ws4.Range("F5").Value = 15346
$ws4.Range("F8").Value = 689
$ws4.Range("F9").Value = 15332
$ws4.Range("F11").Value = 8911
$ws4.Range("F18").Value = 192
$ws4.Range("F34").Value = 45
$ws4.Range("I36").Value = "//i0.hdslb.com/bfs/openplatform/202409/52AMZyUi1727059410434.jpeg"
$ws4.Range("F37").Value = 295
$ws4.Range("F40").Value = 5452
